$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (AD1:AF1), copying the existing header
# cell's style (bold, bordered, centered) so the new headers match the
# rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
for ($i = 2; $i -le 49; $i++) {
    $ws.Range("AD$i").Value = 67
    $ws.Range("AE$i").Value = 95
    $ws.Range("AF$i").Value = 0
}
